# Update automatico via Actualizar 10-18-2020 04-33-14
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# New data rows appended to the "Condicion_Pacientes" table.
$data = @(
    @(44119, 1576, 631, 548, 119, 26),
    @(44120, 1994, 602, 551, 124, 26),
    @(44121, 1950, 903, 534, 128, 23)
)

# Copy the formatting (date style / centered number style) from the last
# existing row down onto the new rows before filling in the values.
$ws.Range("A215:F215").Copy()
$ws.Range("A216:F218").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$r = 216
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $r = $r + 1
}

# Grow the table (and its autofilter) to cover the newly added rows.
$lo.Resize($ws.Range("A1:F218"))

# Reflect the view/selection state captured in the saved workbook.
$ws.Activate() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 198
$win.ScrollColumn = 1
$ws.Range("A219").Select() | Out-Null
